# Adding new CSV utility
# Replace the old sample (Id/Name/Score) sheet contents with a new
# Name/Email/Country/State table imported from a CSV source, and
# right-align every populated cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("John Wick",     "john@test.com",    "US",  "California"),
    @("Sachin Taware",  "Sachin@Test.com", "IND", "Maharashtra")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        $cell.Value = $row[$c]
        $cell.HorizontalAlignment = -4152   # xlRight
    }
}
